$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.644.86"
$ws.Range("E2").Value = "  +1.19%  "

$ws.Range("D3").Value = "1.827.92"
$ws.Range("E3").Value = "  +1.96%  "

$ws.Range("E4").Value = "  +0.32%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "309.16"
$ws.Range("E5").Value = "  +0.82%  "

$ws.Range("E6").Value = "  +0.38%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4684"
$ws.Range("E7").Value = "  +3.65%  "

$ws.Range("E8").Value = "  +0.20%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07147"
$ws.Range("E9").Value = "  +1.14%  "

$ws.Range("E10").Value = "  +2.22%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07675"
$ws.Range("E11").Value = "  -0.94%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "19.43"
$ws.Range("E12").Value = "  -0.07%  "

$ws.Range("D13").Value = "1.872.00"
$ws.Range("E13").Value = "  +3.58%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.258"
$ws.Range("E14").Value = "  -0.41%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.364"
$ws.Range("E15").Value = "  +0.72%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "87.72"
$ws.Range("E16").Value = "  +3.47%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.009"
$ws.Range("E17").Value = "  +0.32%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008562"
$ws.Range("E18").Value = "  +0.83%  "

$ws.Range("E19").Value = "  +0.40%  "

$ws.Range("D20").Value = "26.644.67"
$ws.Range("E20").Value = "  +1.08%  "

$ws.Range("E21").Value = "  +0.08%  "

$ws.Range("E22").Value = "  +1.24%  "

$ws.Range("E23").Value = "  +0.25%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.911"
$ws.Range("E24").Value = "  -2.65%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "152.89"
$ws.Range("E25").Value = "  +1.15%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "17.91"
$ws.Range("E26").Value = "  +0.64%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.000"
$ws.Range("E27").Value = "  -0.94%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "113.73"
$ws.Range("E28").Value = "  +1.72%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.876"
$ws.Range("E29").Value = "  +0.83%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.08819"
$ws.Range("E30").Value = "  +1.59%  "

$ws.Range("E31").Value = "  +1.46%  "

$ws.Range("E32").Value = "  +1.82%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.160"
$ws.Range("E33").Value = "  +5.58%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7361"
$ws.Range("E34").Value = "  +2.14%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.431"
$ws.Range("E35").Value = "  -0.09%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.075"
$ws.Range("E36").Value = "  +0.92%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01931"
$ws.Range("E37").Value = "  +0.28%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.952"
$ws.Range("E38").Value = "  +3.33%  "

$ws.Range("E39").Value = "  +1.49%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.863"
$ws.Range("E40").Value = "  +0.59%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5054"
$ws.Range("E41").Value = "  +0.03%  "

$ws.Range("E42").Value = "  -0.77%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.076"
$ws.Range("E43").Value = "  +0.91%  "

$ws.Range("E44").Value = "  +0.46%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4656"
$ws.Range("E45").Value = "  +0.95%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.06"
$ws.Range("E46").Value = "  +1.09%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "98.42"
$ws.Range("E47").Value = "  -2.08%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.571"
$ws.Range("E48").Value = "  +0.67%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06024"
$ws.Range("E49").Value = "  +1.07%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "63.79"
$ws.Range("E50").Value = "  -0.04%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "35.85"
$ws.Range("E51").Value = "  -0.48%  "
